$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 1020
$ws.Range("J17").Value = 1020
$ws.Range("L17").Value = 3060
$ws.Range("N17").Value = -3396

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 4149.952
$ws.Range("I40").Value = 3049.8
$ws.Range("J40").Value = 4493.75
$ws.Range("K40").Value = 3049.8
$ws.Range("L40").Value = 4493.75
$ws.Range("M40").Value = -2874.8
$ws.Range("N40").Value = -4843.75

# Row 44: Alive and Unwell
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = ""

# Row 87: There Was a Late Fee
$ws.Range("H87").Value = 124990
$ws.Range("J87").Value = 124990
$ws.Range("L87").Value = 124990
$ws.Range("N87").Value = -127486

# Row 90: A Gate Arcane Is Dragon's Bane (L)
$ws.Range("H90").Value = 124990
$ws.Range("J90").Value = 124990
$ws.Range("L90").Value = 374970
$ws.Range("N90").Value = -387450

# Row 127: Liquid Competence
$ws.Range("H127").Value = 1856.9
$ws.Range("I127").Value = 1696.125
$ws.Range("K127").Value = 5088.375
$ws.Range("M127").Value = -128.375

# Row 138: All-night Crafting
$ws.Range("H138").Value = 2201.0977
$ws.Range("J138").Value = 2801.8076
$ws.Range("L138").Value = 8405.4228
$ws.Range("N138").Value = -18685.4228

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 10206313
$ws.Range("I32").Value = 10206313
$ws.Range("K32").Value = 10206313
$ws.Range("M32").Value = -10206026

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 7832122
$ws.Range("I61").Value = 9262594
$ws.Range("J61").Value = 107577.1
$ws.Range("K61").Value = 9262594
$ws.Range("L61").Value = 107577.1
$ws.Range("M61").Value = -9262382
$ws.Range("N61").Value = -108001.1

# Row 63: Rivets Run through It
$ws.Range("H63").Value = 4920.778
$ws.Range("I63").Value = 4410.875
$ws.Range("K63").Value = 4410.875
$ws.Range("M63").Value = -3724.875

# Row 66: A Riveting Revival (L)
$ws.Range("H66").Value = 4920.778
$ws.Range("I66").Value = 4410.875
$ws.Range("K66").Value = 22054.375
$ws.Range("M66").Value = -18622.375

# Row 111: Hedging Bets
$ws.Range("H111").Value = 81496.664
$ws.Range("J111").Value = 81496.664
$ws.Range("L111").Value = 81496.664
$ws.Range("N111").Value = -89676.664

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 5258.486
$ws.Range("I132").Value = 2191.6785
$ws.Range("K132").Value = 6575.0355
$ws.Range("M132").Value = -4045.0355

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 7832122
$ws.Range("I136").Value = 9262594
$ws.Range("J136").Value = 107577.1
$ws.Range("K136").Value = 27787782
$ws.Range("L136").Value = 322731.3
$ws.Range("M136").Value = -27785232
$ws.Range("N136").Value = -327831.3

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 2788.12
$ws.Range("I20").Value = 2595.1428
$ws.Range("J20").Value = 3801.25
$ws.Range("K20").Value = 2595.1428
$ws.Range("L20").Value = 3801.25
$ws.Range("M20").Value = -2348.1428
$ws.Range("N20").Value = -4295.25

# Row 80: Unbreaker
$ws.Range("H80").Value = 1312.5333
$ws.Range("J80").Value = 1361.909
$ws.Range("L80").Value = 1361.909
$ws.Range("N80").Value = -3357.909

# Row 83: Attack on Titanium (L)
$ws.Range("H83").Value = 1312.5333
$ws.Range("J83").Value = 1361.909
$ws.Range("L83").Value = 6809.545
$ws.Range("N83").Value = -16793.545

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 53792.55
$ws.Range("I134").Value = 745.2
$ws.Range("K134").Value = 2235.6
$ws.Range("M134").Value = 299.3999999999996

$ws = $wb.Worksheets.Item("CRP")
# Row 6: Got Your Back
$ws.Range("H6").Value = 20077
$ws.Range("I6").Value = 11877.5
$ws.Range("K6").Value = 11877.5
$ws.Range("M6").Value = -11764.5

# Row 26: As the Worm Turns
$ws.Range("H26").Value = 14476.19
$ws.Range("J26").Value = 14476.19
$ws.Range("L26").Value = 14476.19
$ws.Range("N26").Value = -15050.19

# Row 31: Wall Not Found
$ws.Range("H31").Value = 655215.4399999999
$ws.Range("I31").Value = 8916.666999999999
$ws.Range("K31").Value = 8916.666999999999
$ws.Range("M31").Value = -8621.666999999999

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 655215.4399999999
$ws.Range("I34").Value = 8916.666999999999
$ws.Range("K34").Value = 8916.666999999999
$ws.Range("M34").Value = -8714.666999999999

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 2049.5386
$ws.Range("I58").Value = 2011.7
$ws.Range("J58").Value = 2175.6667
$ws.Range("K58").Value = 2011.7
$ws.Range("L58").Value = 2175.6667
$ws.Range("M58").Value = -1808.7
$ws.Range("N58").Value = -2581.6667

# Row 99: O Pine
$ws.Range("H99").Value = 4108.7856
$ws.Range("I99").Value = 4060.5715
$ws.Range("J99").Value = 4157
$ws.Range("K99").Value = 4060.5715
$ws.Range("L99").Value = 4157
$ws.Range("M99").Value = -2562.5715
$ws.Range("N99").Value = -7153

# Row 126: A Better Conductor
$ws.Range("H126").Value = 4108.7856
$ws.Range("I126").Value = 4060.5715
$ws.Range("J126").Value = 4157
$ws.Range("K126").Value = 12181.7145
$ws.Range("L126").Value = 12471
$ws.Range("M126").Value = -9711.7145
$ws.Range("N126").Value = -17411

# Row 136: Turali Quality
$ws.Range("H136").Value = 2049.5386
$ws.Range("I136").Value = 2011.7
$ws.Range("J136").Value = 2175.6667
$ws.Range("K136").Value = 6035.1
$ws.Range("L136").Value = 6527.000100000001
$ws.Range("M136").Value = -3485.1
$ws.Range("N136").Value = -11627.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service
$ws.Range("H107").Value = 451.84616
$ws.Range("I107").Value = 375.63635
$ws.Range("J107").Value = 507.73334
$ws.Range("K107").Value = 1126.90905
$ws.Range("L107").Value = 1523.20002
$ws.Range("M107").Value = 793.09095
$ws.Range("N107").Value = -5363.20002

# Row 132: More Mezcal
$ws.Range("H132").Value = 1308.2778
$ws.Range("J132").Value = 1409.625
$ws.Range("L132").Value = 12686.625
$ws.Range("N132").Value = -17746.625

$ws = $wb.Worksheets.Item("GSM")
# Row 112: Gentleman Donor
$ws.Range("H112").Value = 15997
$ws.Range("J112").Value = 15997
$ws.Range("L112").Value = 15997
$ws.Range("N112").Value = -18213

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 5206.52
$ws.Range("I46").Value = 3366.0715
$ws.Range("K46").Value = 3366.0715
$ws.Range("M46").Value = -3178.0715

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 37041136
$ws.Range("I55").Value = 52637236
$ws.Range("K55").Value = 52637236
$ws.Range("M55").Value = -52637063

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 733.5714
$ws.Range("I82").Value = 541.55554
$ws.Range("K82").Value = 541.55554
$ws.Range("M82").Value = -180.55554

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 733.5714
$ws.Range("I85").Value = 541.55554
$ws.Range("K85").Value = 541.55554
$ws.Range("M85").Value = 706.44446

# Row 122: Hell on Leather
$ws.Range("H122").Value = 4832.3213
$ws.Range("I122").Value = 4604.7646
$ws.Range("K122").Value = 13814.2938
$ws.Range("M122").Value = -11364.2938

$ws = $wb.Worksheets.Item("WVR")
# Row 7: Long Hair, Long Life
$ws.Range("H7").Value = 17498.75
$ws.Range("I7").Value = 40000
$ws.Range("J7").Value = 9998.333000000001
$ws.Range("K7").Value = 40000
$ws.Range("L7").Value = 9998.333000000001
$ws.Range("M7").Value = -39887
$ws.Range("N7").Value = -10224.333

# Row 40: Our Man in Ul'dah
$ws.Range("H40").Value = 30495
$ws.Range("J40").Value = 30495
$ws.Range("L40").Value = 30495
$ws.Range("N40").Value = -30793
